# Add a "% of Q Drop's" column (I) to the grade-distribution sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + per-course Q-drop percentages (only rows that carry data;
# the section/blank rows in between are left untouched).
$values = @{
    1  = "% of Q Drop's"
    3  = "5.08%"
    6  = "33.33%"
    9  = "3.85%"
    12 = "0.00%"
    15 = "0.00%"
    18 = "0.00%"
    21 = "0.00%"
    24 = "0.00%"
    27 = "0.00%"
    30 = "0.00%"
    31 = "0.00%"
    34 = "0.00%"
    35 = "0.00%"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("I$row")
    # Store as literal text (e.g. "5.08%"), not an auto-converted number.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$row]
}
